$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" (strikeout) values, replacing the old "Strike#" values.
# Regenerated values per commit: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals"
$newK = @{
    2  = 0
    3  = 1
    4  = 5
    5  = 1
    6  = 2
    7  = 5
    8  = 0
    9  = 2
    10 = 3
    11 = 1
    12 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
